$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dayTradeMargins = @{
    2 = "$0.00"
    4 = "$0.00"
    5 = "$0.00"
    6 = "$0.00"
    7 = "$0.00"
    8 = "$0.00"
    9 = "$0.00"
    10 = "$0.00"
    16 = "$990.00"
    17 = "$1,485.00"
    18 = "$852.50"
    19 = "$1,251.50"
    20 = "$2,200.00"
    22 = "$825.00"
    23 = "$2,475.00"
    24 = "$1,375.00"
    25 = "$626.00"
    26 = "$1,100.00"
    30 = "$3,300.00"
    31 = "$6,875.00"
    32 = "$1,100.00"
    33 = "$1,952.50"
    34 = "$467.50"
    35 = "$344.00"
    37 = "$5,032.50"
    38 = "$4,950.00"
    39 = "$2,062.50"
    40 = "$24,200.00"
    41 = "$2,475.00"
}

foreach ($row in $dayTradeMargins.Keys) {
    $cell = $ws.Range("F$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dayTradeMargins[$row]
}

Write-Host "Updated $($dayTradeMargins.Count) DayTradeMargin cells"
